$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: remove C2 value, update E2
$ws.Range("C2").ClearContents()
$ws.Range("E2").Value = 6.253707197847591

# Row 3
$ws.Range("C3").Value = -7.921319741078636
$ws.Range("E3").Value = -14.35806537048446

# Row 5
$ws.Range("E5").Value = 6.136355062499965

# Row 6
$ws.Range("C6").Value = 1.477633171193116

# Row 7
$ws.Range("C7").Value = -0.2005250704869121

# Row 8
$ws.Range("C8").Value = 2.234527904461148

# Row 9
$ws.Range("C9").Value = 1.311727872618218
$ws.Range("E9").Value = 1.093673275363694

# Row 10
$ws.Range("C10").Value = 1.784808447869191

# Row 12
$ws.Range("C12").Value = 2.159589514946725

# Row 13
$ws.Range("C13").Value = 0.8014493436638848
$ws.Range("E13").Value = 1.609625625599986

# Row 14
$ws.Range("C14").Value = -3.107661574595766
$ws.Range("E14").Value = -8.513835774400015

# Row 15
$ws.Range("C15").Value = -1.621578487659103

# Row 16
$ws.Range("C16").Value = 1.906376895025041
$ws.Range("E16").Value = 0.9449384537270955

# Row 17
$ws.Range("C17").Value = -0.4925007786849234
